$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "38.730.80"; E = "  +0.26%  " },
    @{ Row = 3;  D = "2.098.96";  E = "  +0.31%  " },
    @{ Row = 4;  D = $null;       E = "  -0.02%  " },
    @{ Row = 5;  D = "227.60";    E = "  -0.56%  " },
    @{ Row = 6;  D = $null;       E = "  +0.35%  " },
    @{ Row = 7;  D = "62.15";     E = "  +1.43%  " },
    @{ Row = 8;  D = $null;       E = "  -0.01%  " },
    @{ Row = 9;  D = "0.389";     E = "  +1.94%  " },
    @{ Row = 10; D = "0.0838";    E = "  -0.34%  " },
    @{ Row = 11; D = $null;       E = "  -1.21%  " },
    @{ Row = 12; D = "15.66";     E = "  +5.48%  " },
    @{ Row = 13; D = "2.411.08";  E = "  +0.51%  " },
    @{ Row = 14; D = "22.02";     E = "  -1.44%  " },
    @{ Row = 15; D = $null;       E = "  +3.49%  " },
    @{ Row = 16; D = $null;       E = "  +0.86%  " },
    @{ Row = 17; D = "2.099.59";  E = "  +0.13%  " },
    @{ Row = 18; D = "38.697.08"; E = "  +0.37%  " },
    @{ Row = 19; D = "71.65";     E = "  +0.87%  " },
    @{ Row = 20; D = "6.11";      E = "  +0.44%  " },
    @{ Row = 21; D = $null;       E = "  +0.48%  " },
    @{ Row = 22; D = "227.74";    E = "  +0.61%  " },
    @{ Row = 24; D = "2.34";      E = "  -3.68%  " },
    @{ Row = 25; D = $null;       E = "  -0.79%  " },
    @{ Row = 26; D = "9.60";      E = "  +1.68%  " },
    @{ Row = 27; D = "171.84";    E = "  +0.85%  " },
    @{ Row = 28; D = $null;       E = "  +2.27%  " },
    @{ Row = 29; D = $null;       E = "  +3.15%  " },
    @{ Row = 30; D = "19.29";     E = "  +0.77%  " },
    @{ Row = 31; D = $null;       E = "  +8.08%  " },
    @{ Row = 32; D = $null;       E = "  +0.15%  " },
    @{ Row = 33; D = $null;       E = "  +1.02%  " },
    @{ Row = 34; D = $null;       E = "  -1.08%  " },
    @{ Row = 35; D = "7.00";      E = "  +6.86%  " },
    @{ Row = 36; D = $null;       E = "  +1.69%  " },
    @{ Row = 37; D = "2.38";      E = "  -0.18%  " },
    @{ Row = 38; D = "3.53";      E = "  -0.92%  " },
    @{ Row = 39; D = $null;       E = "  +0.17%  " },
    @{ Row = 40; D = $null;       E = "  -2.52%  " },
    @{ Row = 41; D = "102.68";    E = "  +2.57%  " },
    @{ Row = 42; D = "0.0227";    E = "  +2.98%  " },
    @{ Row = 43; D = $null;       E = "  -1.20%  " },
    @{ Row = 44; D = $null;       E = "  +6.70%  " },
    @{ Row = 45; D = $null;       E = "  -0.89%  " },
    @{ Row = 46; D = "7.80";      E = "  +1.75%  " },
    @{ Row = 47; D = "0.0909";    E = "  -0.74%  " },
    @{ Row = 48; D = $null;       E = "  -0.70%  " },
    @{ Row = 49; D = $null;       E = "  +1.59%  " },
    @{ Row = 50; D = $null;       E = "  -0.92%  " },
    @{ Row = 51; D = "2.296.88";  E = "  +0.45%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
